$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.162.77'
$ws.Range('E2').Value = '  +0.45%  '

$ws.Range('D3').Value = '1.801.82'
$ws.Range('E3').Value = '  +2.49%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.08'
$ws.Range('E5').Value = '  +0.61%  '

$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4763'
$ws.Range('E7').Value = '  +26.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3747'
$ws.Range('E8').Value = '  +11.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.42'
$ws.Range('E9').Value = '  +0.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07682'
$ws.Range('E10').Value = '  +6.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.150'
$ws.Range('E11').Value = '  +2.67%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.68'
$ws.Range('E12').Value = '  +0.40%  '

$ws.Range('E13').Value = '  +0.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.401'
$ws.Range('E14').Value = '  +3.83%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.415'
$ws.Range('E15').Value = '  +2.94%  '

$ws.Range('D16').Value = '1.797.79'
$ws.Range('E16').Value = '  +2.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001098'
$ws.Range('E17').Value = '  +3.88%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06755'
$ws.Range('E18').Value = '  +2.39%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.65'
$ws.Range('E19').Value = '  +2.33%  '

$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.52'
$ws.Range('E21').Value = '  +3.22%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.449'
$ws.Range('E22').Value = '  +3.07%  '

$ws.Range('D23').Value = '28.143.59'
$ws.Range('E23').Value = '  +0.34%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.95'
$ws.Range('E24').Value = '  +2.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.410'
$ws.Range('E25').Value = '  +0.38%  '

$ws.Range('E26').Value = '  +5.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.402'
$ws.Range('E27').Value = '  +3.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.43'
$ws.Range('E28').Value = '  -1.22%  '

$ws.Range('D29').Value = '2.004.31'
$ws.Range('E29').Value = '  +2.29%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.63'
$ws.Range('E30').Value = '  +2.08%  '

$ws.Range('E31').Value = '  +0.98%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.049'
$ws.Range('E32').Value = '  +0.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09711'
$ws.Range('E33').Value = '  +10.79%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.959'
$ws.Range('E34').Value = '  +2.71%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02392'
$ws.Range('E35').Value = '  +2.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.24'
$ws.Range('E36').Value = '  +0.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2229'
$ws.Range('E37').Value = '  +5.42%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06380'
$ws.Range('E38').Value = '  +2.83%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.291'
$ws.Range('E39').Value = '  +2.30%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6734'
$ws.Range('E40').Value = '  +1.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.241'
$ws.Range('E41').Value = '  +1.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.483'
$ws.Range('E42').Value = '  +2.53%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.106'
$ws.Range('E43').Value = '  +1.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.24'
$ws.Range('E44').Value = '  +3.55%  '

$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6188'
$ws.Range('E46').Value = '  +2.21%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.864'
$ws.Range('E47').Value = '  +0.89%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.86'
$ws.Range('E48').Value = '  +1.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.069'
$ws.Range('E49').Value = '  +2.70%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.187'
$ws.Range('E50').Value = '  -0.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07130'
$ws.Range('E51').Value = '  -1.08%  '
